$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header row text: drop the old title-cased "Section" / "Question " /
# "Input Type" labels in favor of lower-cased versions ("Question " also had a
# trailing space that gets dropped).
$ws.Range("A1").Value = "section"
$ws.Range("B1").Value = "question"
$ws.Range("C1").Value = "input type"

# Reset the view: select the full A:E columns (instead of the old single-cell
# selection at C41) and scroll back so row 1 / column A is the top-left
# visible cell again (instead of being scrolled down to row 18).
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A1:E1048576").Select() | Out-Null
